# Edit script implementing the diff described in the task.
$d = $word.ActiveDocument

function Replace-Text {
    param(
        [string]$Find,
        [string]$Replace
    )
    $range = $d.Content
    $ok = $range.Find.Execute($Find, $true, $false, $false, $false, $false, $true, 1, $false, $Replace, 2)
    if (-not $ok) {
        Write-Output "WARNING: replace failed for: $Find"
    }
}

# 1. Title "Alice.." -> "Alice."
Replace-Text "Alice.." "Alice."

# 2. laptop or VR goggles -> laptop, or VR goggles
Replace-Text "mobile device, laptop or VR goggles." "mobile device, laptop, or VR goggles."

# 3/4. BBC documentary paragraph: "after a study" -> "after study"; "disregarded it" -> "disregarded them"
Replace-Text "a particular area after a study. It was stated" "a particular area after study. It was stated"
Replace-Text "they disregarded it because they found it uncomfortable" "they disregarded them because they found it uncomfortable"

# 5. "evaluation of a user" -> "evaluations on the user"
Replace-Text "conduct a psychological evaluation of a user." "conduct a psychological evaluations on the user."

# 6. personality traits sentence
Replace-Text "your personality traits." "your personality traits or mental health status."

# 7/8. continous -> continuous ; create a bot -> create a chatbot
Replace-Text "continous" "continuous"
Replace-Text "create a bot with guarantees" "create a chatbot with guarantees"

# 9. human body! -> human body.
Replace-Text "except without a human body!" "except without a human body."

# 10. internet connection / battery life
Replace-Text "The only requirement is internet connection and batery life" "The only requirement is an internet connection and battery life"

# 11. "functions :" -> "functions:"
Replace-Text "functions :" "functions:"

# 12. certain type -> particular type
Replace-Text "asks the chatbot for a certain type of information" "asks the chatbot for a particular type of information"

# 13. entities or keywords that will allow -> entities or keywords will allow
Replace-Text "but instead, entities or keywords that will allow the bot" "but instead, entities or keywords will allow the bot"

# 14. User instructions
Replace-Text "User instructions :" "User instructions:"

# 15a-d. overview paragraph fixes
Replace-Text "an overview about the dialogue" "an overview of the dialogue"
Replace-Text "long text entries, which would make it difficult" "long text entries, making it difficult"
Replace-Text "may be asked for data that he or she may have" "may be asked for pieces that he or she may have"
Replace-Text "in the same text entry the client may report several data (e.g." "in the same text entry, the client may report several pieces of data (e.g."

# 16. first glance
Replace-Text "Thus, from a first glance one can know the subject" "Thus, from first glance, one can know the subject"

# 17a-c. streamlined tool paragraph
Replace-Text "streamlined thanks to the tool, since it will not" "streamlined thanks to the when tool since it will not"
Replace-Text "situation will be avoided where the user starts therapy" "situation will be avoided when the user starts therapy"
Replace-Text "even make a first approach to explaining" "even make the first approach to explaining"

# 18. guiding / interacting -> guide / interact
Replace-Text "guiding the conversation, interacting with language" "guide the conversation, interact with language"

# 19. Alice is sensitive -> Alice might be sensitive
Replace-Text "Alice is sensitive to changes in the wording" "Alice might be sensitive to changes in the wording"

# 20. harmful outputs -> harmful results
Replace-Text "any harmful outputs that could occur" "any harmful results that could occur"

Write-Output "Text replacements complete"
